# Apply changes to the "Customers" sheet: mark additional customers'
# "Test Result" (column F) as TRUE (finished rewriting customer testing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customers")

# Row 2 was previously FALSE -> now TRUE
$ws.Range("F2").Value = $true

# Rows 6 through 14 previously had no value in column F -> now TRUE
$ws.Range("F6").Value = $true
$ws.Range("F7").Value = $true
$ws.Range("F8").Value = $true
$ws.Range("F9").Value = $true
$ws.Range("F10").Value = $true
$ws.Range("F11").Value = $true
$ws.Range("F12").Value = $true
$ws.Range("F13").Value = $true
$ws.Range("F14").Value = $true
